$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 26, pushing the old rows 27-48 down
# to 29-50 (formatting/values carried along automatically by the native
# row insert).
$ws.Rows.Item(27).Insert()
$ws.Rows.Item(27).Insert()

# New row 27: Inferno / Segunda, week of 2021-10-25 (serial 44494)
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C27").Value = "Arica y Parinacota"
$ws.Range("D27").Value = 44494
$ws.Range("D27").NumberFormat = $ws.Range("D29").NumberFormat
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = 100112021
$ws.Range("G27").Value = "Ají"
$ws.Range("H27").Value = "Inferno"
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 120
$ws.Range("K27").Value = 12000
$ws.Range("L27").Value = 13000
$ws.Range("M27").Value = 12500
$ws.Range("N27").Value = "$/caja 15 kilos"
$ws.Range("O27").Value = "Región de Arica y Parinacota"
$ws.Range("P27").Value = 833
$ws.Range("Q27").Value = 15
$ws.Range("R27").Value = "Hortaliza"

# New row 28: Inferno / Tercera, same week (serial 44494)
$ws.Range("A28").Value = 1
$ws.Range("B28").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C28").Value = "Arica y Parinacota"
$ws.Range("D28").Value = 44494
$ws.Range("D28").NumberFormat = $ws.Range("D29").NumberFormat
$ws.Range("E28").Value = 15
$ws.Range("F28").Value = 100112021
$ws.Range("G28").Value = "Ají"
$ws.Range("H28").Value = "Inferno"
$ws.Range("I28").Value = "Tercera"
$ws.Range("J28").Value = 120
$ws.Range("K28").Value = 9000
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = 9500
$ws.Range("N28").Value = "$/caja 15 kilos"
$ws.Range("O28").Value = "Región de Arica y Parinacota"
$ws.Range("P28").Value = 633
$ws.Range("Q28").Value = 15
$ws.Range("R28").Value = "Hortaliza"
